# Update 28-Jan-2021, end of day update.
# Adds the 27-Jan (rows 21-28) and 28-Jan (rows 29-32) petty-cash entries
# to "Sheet1" of the petty cash book, and moves the active selection to C25.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 27-Jan-2021 entries (A21 date already present) ---

# Row 21: Wages Expense
$ws.Range("B21").Value = "Wages Expense"
$ws.Range("D21").Formula = "=60000+260000"

# Row 22: A/R
$ws.Range("B22").Value = "A/R"
$ws.Range("C22").Formula = "=1831000+10000000+57476000"

# Row 23: TRANSFER BCA AA
$ws.Range("B23").Value = "TRANSFER BCA AA"
$ws.Range("D23").Formula = "=1831000"

# Row 24: TRANSFER BCA
$ws.Range("B24").Value = "TRANSFER BCA"
$ws.Range("D24").Formula = "=27770000+1096000"

# Row 25: andreas - prive
$ws.Range("B25").Value = "andreas - prive"
$ws.Range("D25").Value = 5100000

# Row 26: SALES - cash/retail
$ws.Range("B26").Value = "SALES - cash/retail"
$ws.Range("C26").Formula = "=23008975+43182025-57476000"

# Row 27: SELISIH - kurang (new description)
$ws.Range("B27").Value = "SELISIH - kurang"
$ws.Range("D27").Value = 28000

# Row 28: SETOR KE BANK
$ws.Range("B28").Value = "SETOR KE BANK"
$ws.Range("D28").Value = 42000000

# --- 28-Jan-2021 entries (new date) ---

# Row 29: new date + Wages Expense
$ws.Range("A29").Value = 44224
$ws.Range("B29").Value = "Wages Expense"
$ws.Range("D29").Formula = "=60000"

# Row 30: A/R
$ws.Range("B30").Value = "A/R"
$ws.Range("C30").Formula = "=13320000+18450000+40274000+29120000+16368000+14266000"

# Row 31: TRANSFER BCA
$ws.Range("B31").Value = "TRANSFER BCA"
$ws.Range("D31").Formula = "=29370000+14266000"

# Row 32: BELI kresek (new description)
$ws.Range("B32").Value = "BELI kresek"
$ws.Range("D32").Formula = "=99000"

# --- Move the active selection to C25 ---
[void]$ws.Range("C25").Select()
